$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 278, shifting existing rows 278:355 down to 279:356
$ws.Rows("278:278").Insert()

# Populate the newly inserted row 278 with the new price-record data
$ws.Range("A278").Value = 5
$ws.Range("B278").Value = "Macroferia Regional de Talca"
$ws.Range("C278").Value = "Maule"
$ws.Range("D278").Value = 44736
$ws.Range("E278").Value = 7
$ws.Range("F278").Value = 100114013
$ws.Range("G278").Value = "Zanahoria"
$ws.Range("H278").Value = "Sin especificar"
$ws.Range("I278").Value = "Primera"
$ws.Range("J278").Value = 500
$ws.Range("K278").Value = 5500
$ws.Range("L278").Value = 5500
$ws.Range("M278").Value = 5500
$ws.Range("N278").Value = "$/saco 20 kilos"
$ws.Range("O278").Value = "Región de Ñuble"
$ws.Range("P278").Value = 275
$ws.Range("Q278").Value = 20
$ws.Range("R278").Value = "Hortaliza"
